# Implements "--catfile" style edit: abbreviate the criteria headers in row 1
# (F1:I1) and resize the columns so the long free-text caption in column I
# (the actual comment text in I3/I4) stays legible once its header is
# shortened from "Kommentar" to "comm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Abbreviate the criteria captions in the header row -------------------
$ws.Cells.Item(1, 6).Value = "cont"   # F1: Inhalt -> cont
$ws.Cells.Item(1, 7).Value = "styl"   # G1: Stil -> styl
$ws.Cells.Item(1, 8).Value = "spel"   # H1: Rechtschreibung -> spel
$ws.Cells.Item(1, 9).Value = "comm"   # I1: Kommentar -> comm

# --- Resize the affected columns -------------------------------------------
# Column width stored in the sheet is (pixels + 5) / 6 for the Arial 10
# default font (6 px per character). Subtracting the 5/6 character offset
# converts the desired stored width back into the ColumnWidth value that
# needs to be set through the object model.
$ws.Columns.Item(6).ColumnWidth = (5.41 - 5 / 6)    # F
$ws.Columns.Item(7).ColumnWidth = (4.64 - 5 / 6)    # G
$ws.Columns.Item(8).ColumnWidth = (5.18 - 5 / 6)    # H
$ws.Columns.Item(9).ColumnWidth = (47.02 - 5 / 6)   # I (wide, holds full comment text)

# --- Update the view: zoom out a bit and move the active selection --------
$excel.ActiveWindow.Zoom = 180
[void]$ws.Range("J1").Select()
